# Apply the "werkgroep" text and layout adjustments to the ZCBS input form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels (row 2) ---
$ws.Range("O2").Value = "Buurt"
$ws.Range("P2").Value = "Adres"
$ws.Range("U2").Value = "Opmerking"
$ws.Range("W2").Value = "Internet"

# --- Update column widths (O, P, U) to the narrower layout ---
# ColumnWidth is expressed in "characters"; the host rounds to whole
# pixels on save, so the inputs below are chosen to land on the pixel
# that is closest to the authored widths (7.14 / 7.42 / 16.27 chars).
$ws.Columns.Item(15).ColumnWidth = 19/3
$ws.Columns.Item(16).ColumnWidth = 20/3
$ws.Columns.Item(21).ColumnWidth = 15.5

# --- Update the view / selection state ---
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("AC1").Select()
